# "Generate Report for Handback" -- fills in the handback columns on the
# per-locale sheets (zh-cn, de-de) and flips the Overview/status text from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Overview sheet: status text for both locale columns / both rows.
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"

# widen Overview's locale-status columns (E, F) to fit the longer text
$ovw.Range("E1").ColumnWidth = 29.1
$ovw.Range("F1").ColumnWidth = 29.1

# ---------------------------------------------------------------------
# 2) zh-cn sheet: Status text, Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both data rows, plus hyperlinks on the
#    newly-populated "Latest Target File" cells.
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

$zh.Range("I2").Value = "14566da5-554c-4dc9-bc4b-06bb5bfa08fd.md"
$zh.Range("J2").Value = "14566da5-554c-4dc9-bc4b-06bb5bfa08fd.b95fb22e8858c4843917371b7b78f4f7ec9aade4.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-30 21:09:49"

$zh.Range("I3").Value = "9818b088-011a-480c-aba8-3398a251fd19.md"
$zh.Range("J3").Value = "9818b088-011a-480c-aba8-3398a251fd19.b85a6ad70e3cec3e5543c9ba8aadfece84647728.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-30 21:09:49"

$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c9c79b7328fd2b15fbbe02c818f5967e0ad5f069/e2e/14566da5-554c-4dc9-bc4b-06bb5bfa08fd.md", "", "", "14566da5-554c-4dc9-bc4b-06bb5bfa08fd.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c9c79b7328fd2b15fbbe02c818f5967e0ad5f069/e2e/9818b088-011a-480c-aba8-3398a251fd19.md", "", "", "9818b088-011a-480c-aba8-3398a251fd19.md") | Out-Null

$zh.Range("C1").ColumnWidth = 29.1
$zh.Range("I1").ColumnWidth = 39.17
$zh.Range("J1").ColumnWidth = 39.17

# ---------------------------------------------------------------------
# 3) de-de sheet: same shape of edit, different datetime stamp for the
#    "Latest Handback DateTime" column (de-de ran later than zh-cn).
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

$de.Range("I2").Value = "14566da5-554c-4dc9-bc4b-06bb5bfa08fd.md"
$de.Range("J2").Value = "14566da5-554c-4dc9-bc4b-06bb5bfa08fd.b95fb22e8858c4843917371b7b78f4f7ec9aade4.de-de.xlf"
$de.Range("K2").Value = "2016-08-30 21:09:56"

$de.Range("I3").Value = "9818b088-011a-480c-aba8-3398a251fd19.md"
$de.Range("J3").Value = "9818b088-011a-480c-aba8-3398a251fd19.b85a6ad70e3cec3e5543c9ba8aadfece84647728.de-de.xlf"
$de.Range("K3").Value = "2016-08-30 21:09:56"

$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c9c79b7328fd2b15fbbe02c818f5967e0ad5f069/e2e/14566da5-554c-4dc9-bc4b-06bb5bfa08fd.md", "", "", "14566da5-554c-4dc9-bc4b-06bb5bfa08fd.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c9c79b7328fd2b15fbbe02c818f5967e0ad5f069/e2e/9818b088-011a-480c-aba8-3398a251fd19.md", "", "", "9818b088-011a-480c-aba8-3398a251fd19.md") | Out-Null

$de.Range("C1").ColumnWidth = 29.1
$de.Range("I1").ColumnWidth = 39.17
$de.Range("J1").ColumnWidth = 39.17
